{"js": "// Applies the resume-rewrite edit described by the commit:\n// \"improved prompts for skills and projects, changed UI of projects and\n// fixed the glitch issue\"\n//\n// Five body paragraphs get their text replaced (some with embedded line\n// breaks via Word.BreakType.line, matching <w:br/> in the OOXML).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Find the target paragraphs by their current (pre-edit) text so the\n// script is resilient to exact indices, falling back gracefully if the\n// document layout differs slightly.\nconst oldSummary =\n  \"A student with a strong enthusiasm for technology and app development, \" +\n  \"leveraging Computer Science Engineering knowledge to drive innovation \" +\n  \"and growth in the tech industry through the application of academic \" +\n  \"foundations.\";\nconst oldExperience =\n  \"Developed an AI-powered resume web application using AI models, \" +\n  \"training and fine-tuning them to generate high-quality content. \" +\n  \"Optimized model output through prompt refinement, leveraging AI \" +\n  \"capabilities to drive innovative resume development. Demonstrated \" +\n  \"expertise in AI model training and deployment, yielding enhanced \" +\n  \"resume generation with measurable impact through successful model \" +\n  \"output optimization.\";\nconst oldEducation =\n  \"Bachelor of Engineering in Computer Science and Engineering, AVIT, May \" +\n  \"2026, GPA 7.1. Relevant coursework includes computer science and \" +\n  \"engineering fundamentals.\";\nconst oldSkills = \"C#, Node, Python\";\nconst oldProjects =\n  \"Developed a QR scanner and generator web-app with TypeScript and Node, \" +\n  \"enabling secure sharing of encrypted information via PIN protection.\";\n\nlet summaryPara = null;\nlet experiencePara = null;\nlet educationPara = null;\nlet skillsPara = null;\nlet projectsPara = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text === oldSummary) summaryPara = paragraphs.items[i];\n  else if (text === oldExperience) experiencePara = paragraphs.items[i];\n  else if (text === oldEducation) educationPara = paragraphs.items[i];\n  else if (text === oldSkills) skillsPara = paragraphs.items[i];\n  else if (text.indexOf(oldProjects) === 0) projectsPara = paragraphs.items[i];\n}\n\n// --- Summary -------------------------------------------------------------\nif (summaryPara) {\n  summaryPara.insertText(\n    \"Dedicated and people-oriented professional with a background in \" +\n      \"education and experience in supporting and communicating with \" +\n      \"individuals, seeking a role in Human/Social Services or a related \" +\n      \"field. Strong interpersonal skills and a passion for helping \" +\n      \"others are key strengths, poised for growth in a dynamic \" +\n      \"environment focused on social assistance and administration.\",\n    Word.InsertLocation.replace\n  );\n}\n\n// --- Experience ------------------------------------------------------------\nif (experiencePara) {\n  experiencePara.insertText(\n    \"Developed an AI web-app for resume enhancement using LLM models, \" +\n      \"training them to identify required resume features and \" +\n      \"fine-tuning with targeted prompts for content generation, \" +\n      \"resulting in improved content creation efficiency.\",\n    Word.InsertLocation.replace\n  );\n}\n\nawait context.sync();\n\n// --- Education (text, line break, text) -----------------------------------\nif (educationPara) {\n  educationPara.insertText(\n    \"Bachelor of Engineering in Computer Science, AVIT, May 2026, GPA 7.1;\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n  educationPara.insertBreak(Word.BreakType.line, Word.InsertLocation.end);\n  await context.sync();\n  educationPara.insertText(\n    \"Intermediate, Narayana Jr College, June 2022, GPA 8.9\",\n    Word.InsertLocation.end\n  );\n  await context.sync();\n}\n\n// --- Skills (text/line-break sequence, including blank \"double\" breaks) ---\nif (skillsPara) {\n  const skillLines = [\n    \"Good with people, can talk and listen well but sometimes nervous,\",\n    \"Know some billing and petty cash handling but not expert,\",\n    \"Experience running group sessions for mental health and social skills,\",\n    \"Can plan and do social activities but not very creative,\",\n    null, // extra blank line -> double <w:br/>\n    \"Able to multitask but sometimes get overwhelmed,\",\n    null, // extra blank line -> double <w:br/>\n    \"Basic computer skills like Microsoft Word and Excel, \",\n  ];\n\n  skillsPara.insertText(skillLines[0], Word.InsertLocation.replace);\n  await context.sync();\n\n  for (let i = 1; i < skillLines.length; i++) {\n    skillsPara.insertBreak(Word.BreakType.line, Word.InsertLocation.end);\n    await context.sync();\n    const line = skillLines[i];\n    if (line !== null) {\n      skillsPara.insertText(line, Word.InsertLocation.end);\n      await context.sync();\n    }\n  }\n}\n\n// --- Projects (collapse three bulleted runs into a single paragraph) ------\nif (projectsPara) {\n  projectsPara.insertText(\n    \"Developed and implemented a QR scanner and generator, leveraging \" +\n      \"TypeScript and Node.js for the backend. Spearheaded backend \" +\n      \"development, ensuring seamless functionality. Additionally, \" +\n      \"contributed to Prediction Pro, a full-stack application built \" +\n      \"with TypeScript, React, and PostgreSQL, demonstrating expertise \" +\n      \"in modern technologies and collaborative skills through frontend \" +\n      \"and database integration. These projects showcased versatility \" +\n      \"in tech stacks, full-stack development capabilities, and \" +\n      \"effective collaboration.\",\n    Word.InsertLocation.replace\n  );\n}\n\nawait context.sync();\n", "ps1": "# Applies the resume-rewrite edit described by the commit:\n# \"improved prompts for skills and projects, changed UI of projects and\n# fixed the glitch issue\"\n#\n# Five body paragraphs get their text replaced; Education and Skills gain\n# embedded line breaks (InsertBreak -> <w:br/>), Projects collapses three\n# bulleted/broken runs into a single run of prose.\n\n$d = $word.ActiveDocument\n\n$oldSummary = \"A student with a strong enthusiasm for technology and app development, leveraging Computer Science Engineering knowledge to drive innovation and growth in the tech industry through the application of academic foundations.\"\n$oldExperience = \"Developed an AI-powered resume web application using AI models, training and fine-tuning them to generate high-quality content. Optimized model output through prompt refinement, leveraging AI capabilities to drive innovative resume development. Demonstrated expertise in AI model training and deployment, yielding enhanced resume generation with measurable impact through successful model output optimization.\"\n$oldEducation = \"Bachelor of Engineering in Computer Science and Engineering, AVIT, May 2026, GPA 7.1. Relevant coursework includes computer science and engineering fundamentals.\"\n$oldSkills = \"C#, Node, Python\"\n$oldProjectsStart = \"Developed a QR scanner and generator web-app with TypeScript and Node\"\n\n$summaryPara = $null\n$experiencePara = $null\n$educationPara = $null\n$skillsPara = $null\n$projectsPara = $null\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $oldSummary) { $summaryPara = $p }\n    elseif ($t -eq $oldExperience) { $experiencePara = $p }\n    elseif ($t -eq $oldEducation) { $educationPara = $p }\n    elseif ($t -eq $oldSkills) { $skillsPara = $p }\n    elseif ($t.StartsWith($oldProjectsStart)) { $projectsPara = $p }\n}\n\n# --- Summary ---------------------------------------------------------------\nif ($summaryPara) {\n    $summaryPara.Range.Text = \"Dedicated and people-oriented professional with a background in education and experience in supporting and communicating with individuals, seeking a role in Human/Social Services or a related field. Strong interpersonal skills and a passion for helping others are key strengths, poised for growth in a dynamic environment focused on social assistance and administration.\"\n}\n\n# --- Experience --------------------------------------------------------------\nif ($experiencePara) {\n    $experiencePara.Range.Text = \"Developed an AI web-app for resume enhancement using LLM models, training them to identify required resume features and fine-tuning with targeted prompts for content generation, resulting in improved content creation efficiency.\"\n}\n\n# --- Education (text, line break, text) -------------------------------------\nif ($educationPara) {\n    $educationPara.Range.Text = \"Bachelor of Engineering in Computer Science, AVIT, May 2026, GPA 7.1;\"\n    $educationPara.Range.InsertBreak(6)  # wdLineBreak\n    $educationPara.Range.InsertAfter(\"Intermediate, Narayana Jr College, June 2022, GPA 8.9\")\n}\n\n# --- Skills (text/line-break sequence, including blank \"double\" breaks) -----\nif ($skillsPara) {\n    $skillsPara.Range.Text = \"Good with people, can talk and listen well but sometimes nervous,\"\n    $skillsPara.Range.InsertBreak(6)\n    $skillsPara.Range.InsertAfter(\"Know some billing and petty cash handling but not expert,\")\n    $skillsPara.Range.InsertBreak(6)\n    $skillsPara.Range.InsertAfter(\"Experience running group sessions for mental health and social skills,\")\n    $skillsPara.Range.InsertBreak(6)\n    $skillsPara.Range.InsertAfter(\"Can plan and do social activities but not very creative,\")\n    $skillsPara.Range.InsertBreak(6)\n    $skillsPara.Range.InsertBreak(6)\n    $skillsPara.Range.InsertAfter(\"Able to multitask but sometimes get overwhelmed,\")\n    $skillsPara.Range.InsertBreak(6)\n    $skillsPara.Range.InsertBreak(6)\n    $skillsPara.Range.InsertAfter(\"Basic computer skills like Microsoft Word and Excel, \")\n}\n\n# --- Projects (collapse three bulleted runs into a single paragraph) --------\nif ($projectsPara) {\n    $projectsPara.Range.Text = \"Developed and implemented a QR scanner and generator, leveraging TypeScript and Node.js for the backend. Spearheaded backend development, ensuring seamless functionality. Additionally, contributed to Prediction Pro, a full-stack application built with TypeScript, React, and PostgreSQL, demonstrating expertise in modern technologies and collaborative skills through frontend and database integration. These projects showcased versatility in tech stacks, full-stack development capabilities, and effective collaboration.\"\n}\n"}
